$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.068.18"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "3.778.26"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.39%  "
$ws.Range("D7").Value = "3.773.35"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "4.406.78"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "3.772.55"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "68.121.19"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.731"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000138"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.35%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.323"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "448.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("E43").Value = "  -8.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.92%  "
$ws.Range("D46").Value = "2.831.75"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0350"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.17%  "
